$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 445.77777
$ws.Cells.Item(2, 9).Value = 568.2
$ws.Cells.Item(2, 11).Value = 568.2
$ws.Cells.Item(2, 13).Value = -455.2
$ws.Cells.Item(28, 8).Value = 1290.7059
$ws.Cells.Item(28, 10).Value = 1835.7142
$ws.Cells.Item(28, 12).Value = 1835.7142
$ws.Cells.Item(28, 14).Value = -2805.7142
$ws.Cells.Item(43, 8).Value = 101770.3
$ws.Cells.Item(43, 10).Value = 101770.3
$ws.Cells.Item(43, 12).Value = 101770.3
$ws.Cells.Item(43, 14).Value = -101908.3
$ws.Cells.Item(58, 8).Value = 3911
$ws.Cells.Item(58, 10).Value = 18000
$ws.Cells.Item(58, 12).Value = 54000
$ws.Cells.Item(58, 14).Value = -54300
$ws.Cells.Item(86, 8).Value = 3690
$ws.Cells.Item(86, 9).Value = 1974.5
$ws.Cells.Item(86, 10).Value = 5160.4287
$ws.Cells.Item(86, 11).Value = 1974.5
$ws.Cells.Item(86, 12).Value = 5160.4287
$ws.Cells.Item(86, 13).Value = -851.5
$ws.Cells.Item(86, 14).Value = -7406.4287
$ws.Cells.Item(89, 8).Value = 3690
$ws.Cells.Item(89, 9).Value = 1974.5
$ws.Cells.Item(89, 10).Value = 5160.4287
$ws.Cells.Item(89, 11).Value = 9872.5
$ws.Cells.Item(89, 12).Value = 25802.1435
$ws.Cells.Item(89, 13).Value = -4256.5
$ws.Cells.Item(89, 14).Value = -37034.14350000001
$ws.Cells.Item(100, 8).Value = 2231.3076
$ws.Cells.Item(100, 9).Value = 1812
$ws.Cells.Item(100, 10).Value = 3174.75
$ws.Cells.Item(100, 11).Value = 1812
$ws.Cells.Item(100, 12).Value = 3174.75
$ws.Cells.Item(100, 13).Value = -1271
$ws.Cells.Item(100, 14).Value = -4256.75
$ws.Cells.Item(111, 8).Value = 5993.923
$ws.Cells.Item(111, 9).Value = 6292.6
$ws.Cells.Item(111, 10).Value = 4998.3335
$ws.Cells.Item(111, 11).Value = 18877.8
$ws.Cells.Item(111, 12).Value = 14995.0005
$ws.Cells.Item(111, 13).Value = -15810.8
$ws.Cells.Item(111, 14).Value = -21129.0005
$ws.Cells.Item(137, 8).Value = 7345.517
$ws.Cells.Item(137, 9).Value = 4430.095
$ws.Cells.Item(137, 10).Value = 14998.5
$ws.Cells.Item(137, 11).Value = 13290.285
$ws.Cells.Item(137, 12).Value = 44995.5
$ws.Cells.Item(137, 13).Value = -10740.285
$ws.Cells.Item(137, 14).Value = -50095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1382.7778
$ws.Cells.Item(2, 9).Value = 921.4286
$ws.Cells.Item(2, 11).Value = 921.4286
$ws.Cells.Item(2, 13).Value = -808.4286
$ws.Cells.Item(32, 8).Value = 464.86
$ws.Cells.Item(32, 9).Value = 495.0909
$ws.Cells.Item(32, 10).Value = 243.16667
$ws.Cells.Item(32, 11).Value = 495.0909
$ws.Cells.Item(32, 12).Value = 243.16667
$ws.Cells.Item(32, 13).Value = -208.0909
$ws.Cells.Item(32, 14).Value = -817.1666700000001
$ws.Cells.Item(113, 8).Value = 82795.336
$ws.Cells.Item(113, 10).Value = 82795.336
$ws.Cells.Item(113, 12).Value = 82795.336
$ws.Cells.Item(113, 14).Value = -91473.336
$ws.Cells.Item(116, 8).Value = 1382.7778
$ws.Cells.Item(116, 9).Value = 921.4286
$ws.Cells.Item(116, 11).Value = 921.4286
$ws.Cells.Item(116, 13).Value = 1372.5714

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1382.7778
$ws.Cells.Item(3, 9).Value = 921.4286
$ws.Cells.Item(3, 11).Value = 921.4286
$ws.Cells.Item(3, 13).Value = -807.4286
$ws.Cells.Item(6, 8).Value = 29986
$ws.Cells.Item(6, 10).Value = 29986
$ws.Cells.Item(6, 12).Value = 29986
$ws.Cells.Item(6, 14).Value = -30212
$ws.Cells.Item(55, 8).Value = 65281.332
$ws.Cells.Item(55, 10).Value = 65281.332
$ws.Cells.Item(55, 12).Value = 65281.332
$ws.Cells.Item(55, 14).Value = -65827.33199999999
$ws.Cells.Item(96, 8).Value = 55154.332
$ws.Cells.Item(96, 9).Value = 8231.875
$ws.Cells.Item(96, 11).Value = 8231.875
$ws.Cells.Item(96, 13).Value = -5485.875
$ws.Cells.Item(115, 8).Value = 29375
$ws.Cells.Item(115, 10).Value = 29375
$ws.Cells.Item(115, 12).Value = 29375
$ws.Cells.Item(115, 14).Value = -32509

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 7785.5264
$ws.Cells.Item(31, 10).Value = 8284.923000000001
$ws.Cells.Item(31, 12).Value = 8284.923000000001
$ws.Cells.Item(31, 14).Value = -8874.923000000001
$ws.Cells.Item(34, 8).Value = 7785.5264
$ws.Cells.Item(34, 10).Value = 8284.923000000001
$ws.Cells.Item(34, 12).Value = 8284.923000000001
$ws.Cells.Item(34, 14).Value = -8688.923000000001
$ws.Cells.Item(44, 8).Value = 15000
$ws.Cells.Item(44, 10).Value = 15000
$ws.Cells.Item(44, 12).Value = 15000
$ws.Cells.Item(44, 14).Value = -15884
$ws.Cells.Item(57, 8).Value = 39997.332
$ws.Cells.Item(103, 8).Value = 19009.4
$ws.Cells.Item(103, 9).Value = 13761.75
$ws.Cells.Item(103, 11).Value = 13761.75
$ws.Cells.Item(103, 13).Value = -12589.75
$ws.Cells.Item(124, 8).Value = 76661.664
$ws.Cells.Item(124, 10).Value = 76661.664
$ws.Cells.Item(124, 12).Value = 76661.664
$ws.Cells.Item(124, 14).Value = -81571.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 37326.5
$ws.Cells.Item(56, 9).Value = 37326.5
$ws.Cells.Item(56, 11).Value = 37326.5
$ws.Cells.Item(56, 13).Value = -36796.5
$ws.Cells.Item(129, 8).Value = 2118.1
$ws.Cells.Item(129, 9).Value = 546.1667
$ws.Cells.Item(129, 11).Value = 1638.5001
$ws.Cells.Item(129, 13).Value = 3361.4999
$ws.Cells.Item(139, 8).Value = 1505.5454
$ws.Cells.Item(139, 9).Value = 1352.8
$ws.Cells.Item(139, 11).Value = 4058.4
$ws.Cells.Item(139, 13).Value = 1081.6
$ws.Cells.Item(140, 8).Value = 776.8
$ws.Cells.Item(140, 9).Value = 776.8
$ws.Cells.Item(140, 11).Value = 2330.4
$ws.Cells.Item(140, 13).Value = 2849.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 5497.5
$ws.Cells.Item(5, 9).Value = 4997.222
$ws.Cells.Item(5, 10).Value = 10000
$ws.Cells.Item(5, 11).Value = 4997.222
$ws.Cells.Item(5, 12).Value = 10000
$ws.Cells.Item(5, 13).Value = -4885.222
$ws.Cells.Item(5, 14).Value = -10224
$ws.Cells.Item(80, 8).Value = 1233.3334
$ws.Cells.Item(80, 9).Value = 950
$ws.Cells.Item(80, 11).Value = 950
$ws.Cells.Item(80, 13).Value = 48
$ws.Cells.Item(83, 8).Value = 1233.3334
$ws.Cells.Item(83, 9).Value = 950
$ws.Cells.Item(83, 11).Value = 4750
$ws.Cells.Item(83, 13).Value = 242
$ws.Cells.Item(106, 8).Value = 53997
$ws.Cells.Item(106, 10).Value = 53997
$ws.Cells.Item(106, 12).Value = 53997
$ws.Cells.Item(106, 14).Value = -56521

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 9).Value = 1193.3636
$ws.Cells.Item(22, 10).Value = 1216.2632
$ws.Cells.Item(22, 11).Value = 1193.3636
$ws.Cells.Item(22, 12).Value = 1216.2632
$ws.Cells.Item(22, 13).Value = -898.3635999999999
$ws.Cells.Item(22, 14).Value = -1806.2632
$ws.Cells.Item(27, 9).Value = 1193.3636
$ws.Cells.Item(27, 10).Value = 1216.2632
$ws.Cells.Item(27, 11).Value = 1193.3636
$ws.Cells.Item(27, 12).Value = 1216.2632
$ws.Cells.Item(27, 13).Value = -1086.3636
$ws.Cells.Item(27, 14).Value = -1430.2632
$ws.Cells.Item(46, 8).Value = 1400.9524
$ws.Cells.Item(46, 9).Value = 824.7059
$ws.Cells.Item(46, 11).Value = 824.7059
$ws.Cells.Item(46, 13).Value = -636.7059
$ws.Cells.Item(55, 8).Value = 604
$ws.Cells.Item(55, 9).Value = 150
$ws.Cells.Item(55, 10).Value = 664.5333000000001
$ws.Cells.Item(55, 11).Value = 150
$ws.Cells.Item(55, 12).Value = 664.5333000000001
$ws.Cells.Item(55, 13).Value = 23
$ws.Cells.Item(55, 14).Value = -1010.5333
$ws.Cells.Item(98, 8).Value = 69101.42999999999
$ws.Cells.Item(98, 10).Value = 51942
$ws.Cells.Item(98, 12).Value = 51942
$ws.Cells.Item(98, 14).Value = -57932
$ws.Cells.Item(132, 8).Value = 8342.027
$ws.Cells.Item(132, 9).Value = 2910.4333
$ws.Cells.Item(132, 11).Value = 8731.2999
$ws.Cells.Item(132, 13).Value = -6201.2999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 1416
$ws.Cells.Item(107, 9).Value = 1416
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 4248
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -2328
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 128414.125
$ws.Cells.Item(122, 9).Value = 146201.86
$ws.Cells.Item(122, 10).Value = 3900
$ws.Cells.Item(122, 11).Value = 438605.58
$ws.Cells.Item(122, 12).Value = 11700
$ws.Cells.Item(122, 13).Value = -436155.58
$ws.Cells.Item(122, 14).Value = -16600
$ws.Cells.Item(132, 8).Value = 1546.9056
$ws.Cells.Item(132, 9).Value = 1333.5682
$ws.Cells.Item(132, 11).Value = 4000.7046
$ws.Cells.Item(132, 13).Value = -1470.7046
